# Applies the "Platformer system implanted from unity platformer tutorial"
# edit: wraps a batch of loanwords/foreign terms with <w:proofErr> spell
# markers (splitting runs around them), appends a new "2020-02-17" dated
# list item after the platformer list entry, and relocates the stray
# _GoBack bookmark from the final "톨레도" paragraph to the end of the
# "총기 시스템 재설계" list item.

$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Paragraph 4: "...전생한 현역군인으로 플레이하는 이세계물"
#   split trailing run so "이세계물" gets spell-check marks
$d.Paragraphs(4).Range.InsertXML(@"
<w:p $ns><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">총기를 가지고 이세계로 전생한 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>현역군인</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">으로 플레이하는 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>이세계물</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@)

# Paragraph 5: "고강도 이세계 판타지 광물로 인한 높은 TTK구현"
$d.Paragraphs(5).Range.InsertXML(@"
<w:p $ns><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">고강도 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>이세계</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 판타지 광물로 인한 높은 </w:t></w:r><w:r><w:t>TTK</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>구현</w:t></w:r></w:p>
"@)

# Paragraph 7: "...총기모딩을 진행하므로..."
$d.Paragraphs(7).Range.InsertXML(@"
<w:p $ns><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>현실의 부품과 달리</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">현지 대장간에서 한번 개조를 진행하고 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>총기모딩을</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 진행하므로 총기규격과 관계없이 부품사용 가능</w:t></w:r></w:p>
"@)

# Paragraph 10: "기존 프로젝트 파악 및 리팩토링 2020-02-15 ~" (list item)
$d.Paragraphs(10).Range.InsertXML(@"
<w:p $ns><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">기존 프로젝트 파악 및 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>리팩토링</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>2020-02-15 ~</w:t></w:r></w:p>
"@)

# Paragraph 11: "플랫포머 시스템 재설계" (list item) + new "2020-02-17" date run
$d.Paragraphs(11).Range.InsertXML(@"
<w:p $ns><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>플랫포머</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 시스템 재설계</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>2020-02-17</w:t></w:r></w:p>
"@)

# Paragraph 12: "총기 시스템 재설계" (list item) + relocated _GoBack bookmark
$d.Paragraphs(12).Range.InsertXML(@"
<w:p $ns><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>총기 시스템 재설계</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@)

# Paragraph 19: "...단단한 바위로 되어있는 골렘 등"
$d.Paragraphs(19).Range.InsertXML(@"
<w:p $ns><w:r><w:t>7</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>시방향 해변에서 대륙중간까지 이어지는 길고 큰 사막,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">단단한 바위로 되어있는 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>골렘</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 등</w:t></w:r></w:p>
"@)

# Paragraph 20: "5시방향에 망그루브 나무와 늪지대, 비교적 약한 몹과 슬로우를 거는 요상한 기믹 위주"
$d.Paragraphs(20).Range.InsertXML(@"
<w:p $ns><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">5시방향에 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>망그루브</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 나무와 늪지대,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">비교적 약한 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>몹과</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 슬로우를 거는 요상한 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>기믹</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 위주</w:t></w:r></w:p>
"@)

# Paragraph 21: "3시방향에 초원과 대제국, 고유한 금속을 이용한 고티어 기사들이 등장하는 지역"
$d.Paragraphs(21).Range.InsertXML(@"
<w:p $ns><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>3시방향에 초원과 대제국,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">고유한 금속을 이용한 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>고티어</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 기사들이 등장하는 지역</w:t></w:r></w:p>
"@)

# Paragraph 25: "메인메뉴에서는 왼쪽에 플레이어 캐릭터 렌더링, 수첩형태로 오른쪽에 살짝 틀어진 버튼셋"
$d.Paragraphs(25).Range.InsertXML(@"
<w:p $ns><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>메인메뉴에서는</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 왼쪽에 플레이어 캐릭터 렌더링,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">수첩형태로 오른쪽에 살짝 틀어진 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>버튼셋</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@)

# Paragraph 27: "톨레도" - drop the empty pPr/rPr and the (now relocated) bookmark,
# wrap the word itself in spell-check marks
$d.Paragraphs(27).Range.InsertXML(@"
<w:p $ns><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>톨레도</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@)
